# Adds new sample rows (172-181) for dates 43941 (2020-04-20) and 43942 (2020-04-21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows data: row, id, date, provincia, D, E, F, G
# E = $null means leave the cell empty (not present in XML)
$rows = @(
    @{ Row = 172; Id = 170; Date = 43941; Prov = "LAS TUNAS";   D = 31;  E = 0;    F = 0; G = 0 },
    @{ Row = 173; Id = 171; Date = 43941; Prov = "HOLGUÍN";     D = 1;   E = 0;    F = 0; G = 0 },
    @{ Row = 174; Id = 172; Date = 43941; Prov = "GRANMA";      D = 13;  E = 0;    F = 0; G = 1 },
    @{ Row = 175; Id = 173; Date = 43941; Prov = "SANTIAGO";    D = 119; E = 0;    F = 0; G = 1 },
    @{ Row = 176; Id = 174; Date = 43941; Prov = "GUANTÁNAMO";  D = 22;  E = 0;    F = 0; G = 0 },
    @{ Row = 177; Id = 175; Date = 43942; Prov = "LAS TUNAS";   D = 25;  E = $null; F = 0; G = 0 },
    @{ Row = 178; Id = 176; Date = 43942; Prov = "HOLGUÍN";     D = 131; E = $null; F = 1; G = 0 },
    @{ Row = 179; Id = 177; Date = 43942; Prov = "GRANMA";      D = 5;   E = $null; F = 0; G = 0 },
    @{ Row = 180; Id = 178; Date = 43942; Prov = "SANTIAGO";    D = 89;  E = $null; F = 0; G = 0 },
    @{ Row = 181; Id = 179; Date = 43942; Prov = "GUANTÁNAMO";  D = 30;  E = $null; F = 1; G = 0 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Date
    $ws.Cells.Item($row, 2).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
    $ws.Cells.Item($row, 3).Value = $r.Prov
    $ws.Cells.Item($row, 4).Value = $r.D
    if ($null -ne $r.E) {
        $ws.Cells.Item($row, 5).Value = $r.E
    }
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}

# Match the look of column A / B in the prior rows by copying their formats down
$ws.Range("A171").Copy()
$ws.Range("A172:A181").PasteSpecial(-4122)
$ws.Range("B171").Copy()
$ws.Range("B172:B181").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the view to match new data extent: frozen pane top-left cell and active selection
$ws.Application.ActiveWindow.Panes.Item(4).ScrollRow = 165
$ws.Range("D181").Select()
